# Updates cryptos list values (price & 1h volume change) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
# Rows 27/28 (EthereumClassic <-> Stellar) swapped position in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $text) {
    $range = $ws.Range($ref)
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

Set-TextValue 'D2' '34.475.36'
Set-TextValue 'E2' '  -0.39%  '
Set-TextValue 'D3' '1.806.28'
Set-TextValue 'E3' '  +0.65%  '
Set-TextValue 'D5' '228.54'
Set-TextValue 'E5' '  +0.68%  '
Set-TextValue 'E6' '  +4.48%  '
Set-TextValue 'E7' '  +0.04%  '
Set-TextValue 'D8' '34.89'
Set-TextValue 'E8' '  +5.61%  '
Set-TextValue 'E9' '  +1.41%  '
Set-TextValue 'D10' '0.0696'
Set-TextValue 'E10' '  +0.10%  '
Set-TextValue 'D11' '0.0952'
Set-TextValue 'E11' '  +0.18%  '
Set-TextValue 'D12' '2.066.63'
Set-TextValue 'E12' '  +0.62%  '
Set-TextValue 'D13' '11.24'
Set-TextValue 'E13' '  +1.58%  '
Set-TextValue 'D14' '1.804.60'
Set-TextValue 'E14' '  -0.25%  '
Set-TextValue 'D15' '0.645'
Set-TextValue 'E15' '  +1.29%  '
Set-TextValue 'D16' '34.458.62'
Set-TextValue 'E16' '  -0.28%  '
Set-TextValue 'D17' '4.38'
Set-TextValue 'E17' '  +2.17%  '
Set-TextValue 'D18' '69.18'
Set-TextValue 'E18' '  +0.53%  '
Set-TextValue 'D19' '0.0₃0799'
Set-TextValue 'E19' '  -0.35%  '
Set-TextValue 'D20' '245.94'
Set-TextValue 'E20' '  -0.96%  '
Set-TextValue 'D21' '11.55'
Set-TextValue 'E21' '  +2.42%  '
Set-TextValue 'E22' '  +0.13%  '
Set-TextValue 'E23' '  -0.26%  '
Set-TextValue 'D24' '174.77'
Set-TextValue 'E24' '  +5.61%  '
Set-TextValue 'D25' '2.11'
Set-TextValue 'E25' '  +1.92%  '
Set-TextValue 'D26' '7.83'
Set-TextValue 'E26' '  +7.30%  '
Set-TextValue 'B27' 'Stellar'
Set-TextValue 'C27' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D27' '0.120'
Set-TextValue 'E27' '  +2.83%  '
Set-TextValue 'B28' 'EthereumClassic'
Set-TextValue 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D28' '16.81'
Set-TextValue 'E28' '  +1.53%  '
Set-TextValue 'E29' '  -0.24%  '
Set-TextValue 'E30' '  -2.96%  '
Set-TextValue 'E31' '  +0.81%  '
Set-TextValue 'E32' '  +0.65%  '
Set-TextValue 'E33' '  +0.55%  '
Set-TextValue 'D35' '1.397.62'
Set-TextValue 'E35' '  -2.03%  '
Set-TextValue 'E36' '  +1.46%  '
Set-TextValue 'D37' '2.53'
Set-TextValue 'E37' '  -2.57%  '
Set-TextValue 'E38' '  +0.26%  '
Set-TextValue 'E39' '  -1.23%  '
Set-TextValue 'D40' '83.49'
Set-TextValue 'E40' '  -2.28%  '
Set-TextValue 'D41' '2.83'
Set-TextValue 'E41' '  +2.81%  '
Set-TextValue 'D42' '0.950'
Set-TextValue 'E42' '  +1.88%  '
Set-TextValue 'D43' '2.41'
Set-TextValue 'E43' '  -0.51%  '
Set-TextValue 'E44' '  -1.27%  '
Set-TextValue 'D45' '1.12'
Set-TextValue 'E45' '  +3.66%  '
Set-TextValue 'D46' '0.0512'
Set-TextValue 'E46' '  -3.00%  '
Set-TextValue 'E47' '  -2.35%  '
Set-TextValue 'D48' '1.966.69'
Set-TextValue 'E48' '  +0.65%  '
Set-TextValue 'D49' '105.02'
Set-TextValue 'E49' '  -0.90%  '
Set-TextValue 'E50' '  +0.02%  '
Set-TextValue 'E51' '  -0.39%  '
